$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 (2025 data) with the latest metrics
$ws.Range("C6").Value = 400
$ws.Range("D6").Value = 307
$ws.Range("E6").Value = 93
$ws.Range("F6").Value = 67.17724288840262
$ws.Range("G6").Value = 23.25
$ws.Range("H6").Value = 76.75
